$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header values (dropping previously blank H1/I1 cells)
# Write I1 first so the shared-strings table gains "isSourceOf" before
# "isDerivedFrom", matching the target ordering (17=isSourceOf, 18=isDerivedFrom).
$ws.Range("I1").Value = "isSourceOf"
$ws.Range("H1").Value = "isDerivedFrom"

# Move the active selection from H1 to H2
$ws.Range("H2").Select()
